$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "327.87"
Set-TextValue "E2" "6.43%"
Set-TextValue "D3" "39.88"
Set-TextValue "E3" "7.90%"
Set-TextValue "D4" "5.728"
Set-TextValue "E4" "11.56%"
Set-TextValue "E5" "4.23%"
Set-TextValue "D6" "4.564"
Set-TextValue "E6" "3.66%"
Set-TextValue "D7" "8.675"
Set-TextValue "E7" "4.35%"
Set-TextValue "D8" "1.964"
Set-TextValue "E8" "5.52%"
Set-TextValue "D9" "2.998"
Set-TextValue "E9" "1.66%"
Set-TextValue "D10" "0.9417"
Set-TextValue "E10" "1.84%"
Set-TextValue "D11" "0.1291"
Set-TextValue "E11" "15.81%"
Set-TextValue "D12" "0.1989"
Set-TextValue "E12" "6.35%"
Set-TextValue "D13" "0.09158"
Set-TextValue "E13" "4.50%"
Set-TextValue "D14" "0.03506"
Set-TextValue "E14" "6.45%"
Set-TextValue "D15" "0.09625"
Set-TextValue "E15" "0.50%"
Set-TextValue "D16" "0.001315"
Set-TextValue "E16" "-5.22%"
Set-TextValue "D17" "0.006130"
Set-TextValue "E17" "1.81%"
Set-TextValue "D18" "3.369"
Set-TextValue "E18" "-0.67%"
Set-TextValue "D19" "0.3533"
Set-TextValue "E19" "2.36%"
Set-TextValue "D20" "7.587"
Set-TextValue "E20" "19.16%"
Set-TextValue "D21" "0.1410"
Set-TextValue "E21" "9.10%"
Set-TextValue "D22" "0.2425"
Set-TextValue "E22" "2.13%"
Set-TextValue "D23" "0.04438"
Set-TextValue "E23" "2.10%"
Set-TextValue "D24" "0.001251"
Set-TextValue "E24" "4.02%"
Set-TextValue "E25" "1.22%"
Set-TextValue "D26" "0.0001190"
Set-TextValue "E26" "-15.33%"
Set-TextValue "D27" "0.0003989"
Set-TextValue "E27" "37.26%"
Set-TextValue "D39" "0.02527"
Set-TextValue "E39" "17.86%"
Set-TextValue "E40" "6.08%"
Set-TextValue "D41" "0.007309"
Set-TextValue "E41" "-3.63%"
Set-TextValue "E42" "5.59%"
Set-TextValue "D43" "0.008884"
Set-TextValue "E43" "4.43%"
Set-TextValue "D44" "0.002189"
Set-TextValue "E44" "9.92%"
Set-TextValue "D45" "0.009988"
Set-TextValue "E45" "15.64%"
Set-TextValue "D46" "0.00006704"
Set-TextValue "E46" "1.93%"
Set-TextValue "D47" "0.00000000750"
Set-TextValue "E47" "-0.39%"
Set-TextValue "D48" "0.002873"
Set-TextValue "E48" "-13.06%"
Set-TextValue "D49" "0.001800"
Set-TextValue "E49" "24.46%"
Set-TextValue "D50" "0.00002099"
Set-TextValue "E50" "-0.39%"
Set-TextValue "D51" "0.0001999"
Set-TextValue "E51" "-0.39%"
